$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-increment the "Parent Public ID" (column A) per statement group,
# replacing the old "Statement 1" / "Statement 6" labels with new
# machine-friendly ids.
$ws.Range("A2:A4").Value = "statement-01"
$ws.Range("A5:A6").Value = "statement-02"

# "Public ID" column (B) is no longer populated for these sample rows.
$ws.Range("B2:B6").ClearContents()
